# Adds a new worksheet ("instrument_calib") holding the analog-transfer-
# function coefficients, placed right after the existing "instrument_info"
# sheet, and makes it the active tab - matching the target commit
# "added sheet for analog transfer function".

$wb = $excel.ActiveWorkbook

# --- existing sheet: move the selection that was left on D18 -----------
$wsInfo = $wb.Worksheets.Item(1)
$wsInfo.Range("A7").Select()

# --- new sheet, inserted after "instrument_info" ------------------------
$wsCalib = $wb.Worksheets.Add($null, $wsInfo)
$wsCalib.Name = "instrument_calib"

# Column widths close to the authored best-fit sizing for the two columns.
$wsCalib.Columns.Item(1).ColumnWidth = 17.3
$wsCalib.Columns.Item(2).ColumnWidth = 20.3

# Populate in the same left-to-right, top-to-bottom order the author used
# (this also reproduces the shared-string table insertion order).
$wsCalib.Range("B1").Value = "Coefficients"
$wsCalib.Range("A1").Value = "Instrument name"
$wsCalib.Range("A2").Value = "coil current"
$wsCalib.Range("B2").Value = "[0,1]"

# Leave the selection/active tab where the author left it when saving.
$wsCalib.Range("B3").Select()
